# Laporan, Rabu 19 Februari 2025
# - Update the "Isi Laporan" cell (C4) text: swap article #1 and #4 references
#   for the new ones used in this week's entry.
# - Row 4 is a bit shorter now that the text is slightly different (210 -> 180).
# - Move the active selection from F4 to D4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lines = @(
    "Mencari 5 artikel yang berkaitan dengan pemodelan matematika dengan tema tentang laju pertumbuhan diantaranya : ",
    "1. Proyeksi Pertumbuhan Penduduk Sulawesi Tenggara dengan menggunakan model eksponensial dan model logistik (https://jaf.uho.ac.id/index.php/JAFUHO/article/view/15/6), ",
    "2. Pengaruh Jumlah Penduduk Miskin (https://journals.unisba.ac.id/index.php/JRIEB/article/view/401), ",
    "3. Pengaruh Upah Minimum (https://journals.unisba.ac.id/index.php/JRIEB/article/view/1911), ",
    "4. Proyeksi Pertumbuhan Mobil Pribadi Roda Empat (Plat Hitam) Kota Manado Menggunakan Persamaan Differensial Model Pertumbuhan Populasi Kontinu (Model Logistik) (https://ejournal.unsrat.ac.id/v3/index.php/decartesian/article/view/14017/13590), ",
    "5. Model Pertumbuhan Populasi Malthus (https://media.neliti.com/media/publications/185154-ID-kestabilan-populasi-model-lotka-volterra.pdf)."
)
$newText = [string]::Join("`n", $lines)

$ws.Range("C4").Value = $newText

# The row shrinks a little to fit the revised text.
$ws.Rows.Item(4).RowHeight = 180

# Active cell / selection moves from F4 to D4.
$null = $ws.Range("D4").Select()
